$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

$ws.Range("C3").Value = 36
$ws.Range("C4").Value = 439
$ws.Range("C5").Value = 27
$ws.Range("C6").Value = 6
$ws.Range("C7").Value = 2073

[void]$ws.Columns.Item(9).Select()
